# feat: added serial number to offline tool
# - Removed CADID column (content renamed to Serial Number, kept in place of the
#   old "Comments" slot) because it wasn't used and was hidden
# - Added serial number column
# - Changed layout so "Comments" moves to the last column and an (initially
#   empty) styled cell is added below it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CAD ID" header text becomes "Serial Number", and that column now
# lives where "Comments" used to be (AD); "Comments" shifts out to the new
# last column (AE).
$ws.Range("AD1").Value = "Serial Number"
$ws.Range("AE1").Value = "Comments"

# Column AE grew to fit the "Comments" header/content that moved into it.
$ws.Columns.Item(31).ColumnWidth = 15.022135416666666

# New empty, underlined cell added at AE4 (new style: underline font, no bold).
$ws.Range("AE4").Font.Underline = $true
$ws.Range("AE4").Value = ""

# Update the active selection to match where the user ended up editing.
$ws.Range("AE4").Select() | Out-Null
